# Apply the updated crypto price/volume data to the sheet (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "23.009.14"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.69%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.601.86"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.83%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.002"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "301.15"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3779"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.87%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3627"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -5.50%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "49.43"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.03%  "
$ws.Range("E10").Value = "  -6.30%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.001"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.11%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08127"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.74%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.83"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.24%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.592"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.10%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.391"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -6.56%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001243"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.49%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.599.53"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.92%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "92.15"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.87%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06876"
$ws.Range("D19").Style = "Normal"
$ws.Range("E20").Value = "  -6.62%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.567"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.35%  "
$ws.Range("B22").Value = "BitDAO"
$ws.Range("C22").Value = "https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.5572"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.66%  "
$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.002"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.16%  "
$ws.Range("B24").Value = "Cosmos"
$ws.Range("C24").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.16"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.51%  "
$ws.Range("B25").Value = "WrappedBTC"
$ws.Range("C25").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "23.003.49"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.69%  "
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.359"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.42%  "
$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.798"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.91%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "21.08"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.01%  "
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "150.33"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.53%  "
$ws.Range("B30").Value = "HuobiToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.253"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.40%  "
$ws.Range("B31").Value = "BitcoinCash"
$ws.Range("C31").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "133.61"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.57%  "
$ws.Range("B32").Value = "WEMIXTOKEN"
$ws.Range("C32").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.314"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.87%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.801"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -11.93%  "
$ws.Range("B34").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C34").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.778.37"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.77%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9613"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.01%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.07631"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.02%  "
$ws.Range("B37").Value = "FraxShare"
$ws.Range("C37").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.39"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.57%  "
$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.292"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.95%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02705"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -7.44%  "
$ws.Range("B40").Value = "Algorand"
$ws.Range("C40").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2537"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.44%  "
$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.08859"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.87%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.364"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.98%  "
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7056"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.51%  "
$ws.Range("B44").Value = "Aptos"
$ws.Range("C44").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "12.49"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.89%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "15.18"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -9.81%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6615"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.59%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.315"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.05%  "
$ws.Range("B48").Value = "Frax"
$ws.Range("C48").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.000"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.02%  "
$ws.Range("B49").Value = "PancakeSwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.992"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.59%  "
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "132.77"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.16%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07904"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.40%  "
